$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# ---------------------------------------------------------------------------
# Cell (1,1): "<< Client Name >>" - single run, whole-paragraph resize 26->24
# ---------------------------------------------------------------------------
$cell11 = $t.Cell(1, 1)
$r11 = $cell11.Range
$r11.Font.Size = 12
$r11.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# Cell (1,2): spaces + "Sneha Shukla" -> split into 2 runs, 2nd one resized
# ---------------------------------------------------------------------------
$cell12 = $t.Cell(1, 2)
$full12 = $cell12.Range
$rng12 = $d.Range($full12.Start, $full12.End - 1)

$xml12 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="03139B0E" w14:textId="77777777" w:rsidR="00B16E98" w:rsidRDefault="00B16E98" w:rsidP="00B16E98">
<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">                                 </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Sneha Shukla</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng12.InsertXML($xml12)

# ---------------------------------------------------------------------------
# Cell (2,1): "<< Date >>" - single run, whole-paragraph resize 26->24
# ---------------------------------------------------------------------------
$cell21 = $t.Cell(2, 1)
$r21 = $cell21.Range
$r21.Font.Size = 12
$r21.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# Cell (2,2): 3 runs (spaces, spaces, "<< Date >>") -> merge into 1 run,
# resized, with combined whitespace trimmed from 40 to 30 spaces.
# ---------------------------------------------------------------------------
$cell22 = $t.Cell(2, 2)
$full22 = $cell22.Range
$rng22 = $d.Range($full22.Start, $full22.End - 1)

$xml22 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="3208FC08" w14:textId="19FC0B83" w:rsidR="00B16E98" w:rsidRDefault="00B16E98" w:rsidP="00B16E98">
<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">                              &lt;&lt; Date &gt;&gt;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng22.InsertXML($xml22)

Write-Output "All edits applied"
